# Apply the "feat: add 2022-Q4 data" edit:
#  1. Insert the new 2022-Q4 summary row at the top of the "总计" (total)
#     sheet's data (columns B:D only -- the leading index column A keeps
#     counting up 0,1,2,... and simply grows by one row at the bottom).
#  2. Insert a brand-new worksheet named "2022-Q4" right after "总计"
#     holding the per-fund breakdown for that quarter.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. "总计" sheet: push existing quarters down one row, add 2022-Q4 on top ---

# New last row (index 7) takes what used to be row 8's data (2020-Q4)
$ws1.Cells.Item(9, 1).Value = 7
$ws1.Range("A8").Copy()
$ws1.Range("A9").PasteSpecial(-4122)
$ws1.Cells.Item(9, 2).Value = "2020-Q4"
$ws1.Cells.Item(9, 3).Value = 10
$ws1.Cells.Item(9, 4).Value = 0.29

$ws1.Cells.Item(8, 2).Value = "2021-Q1"
$ws1.Cells.Item(8, 3).Value = 7
$ws1.Cells.Item(8, 4).Value = 0.47

$ws1.Cells.Item(7, 2).Value = "2021-Q2"
$ws1.Cells.Item(7, 3).Value = 2
$ws1.Cells.Item(7, 4).Value = 0.04

$ws1.Cells.Item(6, 2).Value = "2021-Q3"
$ws1.Cells.Item(6, 3).Value = 4
$ws1.Cells.Item(6, 4).Value = 0.02

$ws1.Cells.Item(5, 2).Value = "2021-Q4"
$ws1.Cells.Item(5, 3).Value = 12
$ws1.Cells.Item(5, 4).Value = 2.32

$ws1.Cells.Item(4, 2).Value = "2022-Q1"
$ws1.Cells.Item(4, 3).Value = 25
$ws1.Cells.Item(4, 4).Value = 4.13

$ws1.Cells.Item(3, 2).Value = "2022-Q2"
$ws1.Cells.Item(3, 3).Value = 24
$ws1.Cells.Item(3, 4).Value = 2.96

$ws1.Cells.Item(2, 2).Value = "2022-Q4"
$ws1.Cells.Item(2, 3).Value = 4
$ws1.Cells.Item(2, 4).Value = 0.39

# --- 2. Insert the new "2022-Q4" worksheet right after "总计" ---

$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q4"

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Numeric-looking columns (B, D, E, F, G) must stay text, matching the
# source data export (inlineStr in the OOXML) rather than being coerced
# to numbers -- this also protects fund codes with leading zeros.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "502023"
$newSheet.Cells.Item(2, 3).Value = "鹏华国证钢铁行业指数（LOF）A"
$newSheet.Cells.Item(2, 4).Value = "8.97"
$newSheet.Cells.Item(2, 5).Value = "94.55"
$newSheet.Cells.Item(2, 6).Value = "2.53"
$newSheet.Cells.Item(2, 7).Value = "0.2269"
$newSheet.Cells.Item(2, 8).Value = 10

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "012810"
$newSheet.Cells.Item(3, 3).Value = "鹏华国证钢铁行业指数（LOF）C"
$newSheet.Cells.Item(3, 4).Value = "3.28"
$newSheet.Cells.Item(3, 5).Value = "94.55"
$newSheet.Cells.Item(3, 6).Value = "2.53"
$newSheet.Cells.Item(3, 7).Value = "0.0830"
$newSheet.Cells.Item(3, 8).Value = 10

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "168203"
$newSheet.Cells.Item(4, 3).Value = "中融国证钢铁行业指数A"
$newSheet.Cells.Item(4, 4).Value = "3.18"
$newSheet.Cells.Item(4, 5).Value = "92.38"
$newSheet.Cells.Item(4, 6).Value = "2.46"
$newSheet.Cells.Item(4, 7).Value = "0.0782"
$newSheet.Cells.Item(4, 8).Value = 10

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "016815"
$newSheet.Cells.Item(5, 3).Value = "中融国证钢铁行业指数C"
$newSheet.Cells.Item(5, 4).Value = "0.02"
$newSheet.Cells.Item(5, 5).Value = "92.38"
$newSheet.Cells.Item(5, 6).Value = "2.46"
$newSheet.Cells.Item(5, 7).Value = "0.0005"
$newSheet.Cells.Item(5, 8).Value = 10

# Match the header row / index-column styling used throughout the workbook
# (bold, centered, bordered cell style applied to B1:H1 and A2:A5).
$ws1.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Restore "总计" as the active sheet/selection (unaffected by this edit).
$null = $newSheet.Range("A1").Select()
$ws1.Activate()
